# Auto-generated edit script applying the Hyperion_Profits market-price refresh diff.
# For each affected (sheet, row) the H..N "leve profit" columns are updated to the
# refreshed currentAveragePrice-derived figures; a few rows also gain/lose an N (or
# M/N) cell entirely, mirrored here with ClearContents() for removals.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 955.7619
$ws.Range("I2").Value = 1083.7858
$ws.Range("J2").Value = 699.7143
$ws.Range("K2").Value = 1083.7858
$ws.Range("L2").Value = 699.7143
$ws.Range("M2").Value = -970.7858000000001
$ws.Range("N2").Value = -925.7143
$ws.Range("H6").Value = 90
$ws.Range("I6").Value = 90
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 270
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -158
$ws.Range("N6").ClearContents()
$ws.Range("H38").Value = 42.75
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H58").Value = 5614.067
$ws.Range("J58").Value = 6986.3335
$ws.Range("L58").Value = 20959.0005
$ws.Range("N58").Value = -21259.0005
$ws.Range("H107").Value = 111112780
$ws.Range("I107").Value = 111112780
$ws.Range("K107").Value = 111112780
$ws.Range("M107").Value = -111110860
$ws.Range("H138").Value = 3079.848
$ws.Range("I138").Value = 1413.4783
$ws.Range("K138").Value = 4240.4349
$ws.Range("M138").Value = 899.5650999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3357.7031
$ws.Range("I32").Value = 1998.3846
$ws.Range("K32").Value = 1998.3846
$ws.Range("M32").Value = -1711.3846
$ws.Range("H45").Value = 9595269
$ws.Range("J45").Value = 9581
$ws.Range("L45").Value = 9581
$ws.Range("N45").Value = -10335
$ws.Range("H61").Value = 3551
$ws.Range("I61").Value = 3018.5715
$ws.Range("J61").Value = 4296.4
$ws.Range("K61").Value = 3018.5715
$ws.Range("L61").Value = 4296.4
$ws.Range("M61").Value = -2806.5715
$ws.Range("N61").Value = -4720.4
$ws.Range("H63").Value = 2431.1667
$ws.Range("I63").Value = 2431.1667
$ws.Range("K63").Value = 2431.1667
$ws.Range("M63").Value = -1745.1667
$ws.Range("H66").Value = 2431.1667
$ws.Range("I66").Value = 2431.1667
$ws.Range("K66").Value = 12155.8335
$ws.Range("M66").Value = -8723.8335
$ws.Range("H88").Value = 1599.1666
$ws.Range("J88").Value = 2170.8572
$ws.Range("L88").Value = 2170.8572
$ws.Range("N88").Value = -2982.8572
$ws.Range("H91").Value = 1599.1666
$ws.Range("J91").Value = 2170.8572
$ws.Range("L91").Value = 2170.8572
$ws.Range("N91").Value = -4978.8572
$ws.Range("H132").Value = 3367.8096
$ws.Range("J132").Value = 4668.3335
$ws.Range("L132").Value = 14005.0005
$ws.Range("N132").Value = -19065.0005
$ws.Range("H136").Value = 3551
$ws.Range("I136").Value = 3018.5715
$ws.Range("J136").Value = 4296.4
$ws.Range("K136").Value = 9055.7145
$ws.Range("L136").Value = 12889.2
$ws.Range("M136").Value = -6505.7145
$ws.Range("N136").Value = -17989.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3971803
$ws.Range("I107").Value = 5954280
$ws.Range("J107").Value = 6849.8335
$ws.Range("K107").Value = 5954280
$ws.Range("L107").Value = 6849.8335
$ws.Range("M107").Value = -5952360
$ws.Range("N107").Value = -10689.8335
$ws.Range("H134").Value = 3963.7856
$ws.Range("I134").Value = 1666.2858
$ws.Range("K134").Value = 4998.857400000001
$ws.Range("M134").Value = -2463.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 15981.146
$ws.Range("I31").Value = 1921.579
$ws.Range("K31").Value = 1921.579
$ws.Range("M31").Value = -1626.579
$ws.Range("H34").Value = 15981.146
$ws.Range("I34").Value = 1921.579
$ws.Range("K34").Value = 1921.579
$ws.Range("M34").Value = -1719.579
$ws.Range("H57").Value = 19250
$ws.Range("J57").Value = 18600
$ws.Range("L57").Value = 18600
$ws.Range("N57").Value = -19720
$ws.Range("H86").Value = 13118.435
$ws.Range("I86").Value = 12544.889
$ws.Range("J86").Value = 13487.143
$ws.Range("K86").Value = 12544.889
$ws.Range("L86").Value = 13487.143
$ws.Range("M86").Value = -11421.889
$ws.Range("N86").Value = -15733.143
$ws.Range("H89").Value = 13118.435
$ws.Range("I89").Value = 12544.889
$ws.Range("J89").Value = 13487.143
$ws.Range("K89").Value = 62724.44499999999
$ws.Range("L89").Value = 67435.715
$ws.Range("M89").Value = -57108.44499999999
$ws.Range("N89").Value = -78667.715
$ws.Range("H107").Value = 1061.7097
$ws.Range("I107").Value = 1023.48
$ws.Range("J107").Value = 1221
$ws.Range("K107").Value = 1023.48
$ws.Range("L107").Value = 1221
$ws.Range("M107").Value = 896.52
$ws.Range("N107").Value = -5061
$ws.Range("H132").Value = 65284.832
$ws.Range("I132").Value = 49888.617
$ws.Range("K132").Value = 149665.851
$ws.Range("M132").Value = -147135.851

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 2813.4443
$ws.Range("I7").Value = 3231.5715
$ws.Range("J7").Value = 1350
$ws.Range("K7").Value = 9694.7145
$ws.Range("L7").Value = 4050
$ws.Range("M7").Value = -9582.7145
$ws.Range("N7").Value = -4274
$ws.Range("H60").Value = 863
$ws.Range("I60").Value = 297
$ws.Range("K60").Value = 891
$ws.Range("M60").Value = -640
$ws.Range("H132").Value = 1659.3334
$ws.Range("I132").Value = 1206.1111
$ws.Range("J132").Value = 2339.1667
$ws.Range("K132").Value = 10854.9999
$ws.Range("L132").Value = 21052.5003
$ws.Range("M132").Value = -8324.9999
$ws.Range("N132").Value = -26112.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9096038
$ws.Range("I70").Value = 11115935
$ws.Range("J70").Value = 6499.75
$ws.Range("K70").Value = 11115935
$ws.Range("L70").Value = 6499.75
$ws.Range("M70").Value = -11115665
$ws.Range("N70").Value = -7039.75
$ws.Range("H73").Value = 9096038
$ws.Range("I73").Value = 11115935
$ws.Range("J73").Value = 6499.75
$ws.Range("K73").Value = 11115935
$ws.Range("L73").Value = 6499.75
$ws.Range("M73").Value = -11114999
$ws.Range("N73").Value = -8371.75
$ws.Range("H80").Value = 1879389.6
$ws.Range("I80").Value = 2711773.2
$ws.Range("K80").Value = 2711773.2
$ws.Range("M80").Value = -2710775.2
$ws.Range("H83").Value = 1879389.6
$ws.Range("I83").Value = 2711773.2
$ws.Range("K83").Value = 13558866
$ws.Range("M83").Value = -13553874
$ws.Range("H132").Value = 3972.6924
$ws.Range("I132").Value = 3468
$ws.Range("J132").Value = 6748.5
$ws.Range("K132").Value = 10404
$ws.Range("L132").Value = 20245.5
$ws.Range("M132").Value = -7874
$ws.Range("N132").Value = -25305.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5042.3
$ws.Range("I7").Value = 1904.1666
$ws.Range("J7").Value = 9749.5
$ws.Range("K7").Value = 1904.1666
$ws.Range("L7").Value = 9749.5
$ws.Range("M7").Value = -1792.1666
$ws.Range("N7").Value = -9973.5
$ws.Range("H68").Value = 1524.8182
$ws.Range("I68").Value = 1829
$ws.Range("J68").Value = 1159.8
$ws.Range("K68").Value = 1829
$ws.Range("L68").Value = 1159.8
$ws.Range("M68").Value = -1080
$ws.Range("N68").Value = -2657.8
$ws.Range("H71").Value = 1524.8182
$ws.Range("I71").Value = 1829
$ws.Range("J71").Value = 1159.8
$ws.Range("K71").Value = 9145
$ws.Range("L71").Value = 5799
$ws.Range("M71").Value = -5401
$ws.Range("N71").Value = -13287
$ws.Range("H126").Value = 5042.3
$ws.Range("I126").Value = 1904.1666
$ws.Range("J126").Value = 9749.5
$ws.Range("K126").Value = 5712.4998
$ws.Range("L126").Value = 29248.5
$ws.Range("M126").Value = -3242.4998
$ws.Range("N126").Value = -34188.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 59899
$ws.Range("J95").Value = 59899
$ws.Range("L95").Value = 59899
$ws.Range("N95").Value = -65391
$ws.Range("H122").Value = 1670.1482
$ws.Range("I122").Value = 1265.4348
$ws.Range("J122").Value = 3997.25
$ws.Range("K122").Value = 3796.3044
$ws.Range("L122").Value = 11991.75
$ws.Range("M122").Value = -1346.3044
$ws.Range("N122").Value = -16891.75
$ws.Range("H132").Value = 63190588
$ws.Range("I132").Value = 111125730
$ws.Range("J132").Value = 1559697.2
$ws.Range("K132").Value = 333377190
$ws.Range("L132").Value = 4679091.6
$ws.Range("M132").Value = -333374660
$ws.Range("N132").Value = -4684151.6
$ws.Range("H136").Value = 4212.143
$ws.Range("I136").Value = 2413.3333
$ws.Range("K136").Value = 7239.999899999999
$ws.Range("M136").Value = -4689.999899999999
